$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Step 1: delete row 10 (the row with 19+34=, 59+3=, 82-8=, 84-39=, 66+7=)
$t.Rows.Item(10).Delete()

# Step 2: update remaining cell equations (old -> new), topologically ordered
# to avoid cross-contamination when a new value equals another cell's old value
$d.Content.Find.Execute("88+3=", $true, $false, $false, $false, $false, $true, 1, $false, "93-38=", 2) | Out-Null
$d.Content.Find.Execute("25+8=", $true, $false, $false, $false, $false, $true, 1, $false, "46-7=", 2) | Out-Null
$d.Content.Find.Execute("14+7=", $true, $false, $false, $false, $false, $true, 1, $false, "80-69=", 2) | Out-Null
$d.Content.Find.Execute("59+35=", $true, $false, $false, $false, $false, $true, 1, $false, "55-19=", 2) | Out-Null
$d.Content.Find.Execute("43-9=", $true, $false, $false, $false, $false, $true, 1, $false, "5+8=", 2) | Out-Null
$d.Content.Find.Execute("5+18=", $true, $false, $false, $false, $false, $true, 1, $false, "79+18=", 2) | Out-Null
$d.Content.Find.Execute("84-59=", $true, $false, $false, $false, $false, $true, 1, $false, "83-66=", 2) | Out-Null
$d.Content.Find.Execute("82-26=", $true, $false, $false, $false, $false, $true, 1, $false, "84-59=", 2) | Out-Null
$d.Content.Find.Execute("59+22=", $true, $false, $false, $false, $false, $true, 1, $false, "95-19=", 2) | Out-Null
$d.Content.Find.Execute("95-7=", $true, $false, $false, $false, $false, $true, 1, $false, "35+26=", 2) | Out-Null
$d.Content.Find.Execute("8+5=", $true, $false, $false, $false, $false, $true, 1, $false, "26+36=", 2) | Out-Null
$d.Content.Find.Execute("24+9=", $true, $false, $false, $false, $false, $true, 1, $false, "70-67=", 2) | Out-Null
$d.Content.Find.Execute("48+36=", $true, $false, $false, $false, $false, $true, 1, $false, "92-85=", 2) | Out-Null
$d.Content.Find.Execute("83-76=", $true, $false, $false, $false, $false, $true, 1, $false, "97-58=", 2) | Out-Null
$d.Content.Find.Execute("14+58=", $true, $false, $false, $false, $false, $true, 1, $false, "92-37=", 2) | Out-Null
$d.Content.Find.Execute("60-55=", $true, $false, $false, $false, $false, $true, 1, $false, "20-8=", 2) | Out-Null
$d.Content.Find.Execute("73+8=", $true, $false, $false, $false, $false, $true, 1, $false, "47-29=", 2) | Out-Null
$d.Content.Find.Execute("38+26=", $true, $false, $false, $false, $false, $true, 1, $false, "90-81=", 2) | Out-Null
$d.Content.Find.Execute("70-36=", $true, $false, $false, $false, $false, $true, 1, $false, "55+7=", 2) | Out-Null
$d.Content.Find.Execute("74-56=", $true, $false, $false, $false, $false, $true, 1, $false, "97-88=", 2) | Out-Null
$d.Content.Find.Execute("38+29=", $true, $false, $false, $false, $false, $true, 1, $false, "70-14=", 2) | Out-Null
$d.Content.Find.Execute("83-58=", $true, $false, $false, $false, $false, $true, 1, $false, "45+38=", 2) | Out-Null
$d.Content.Find.Execute("51-36=", $true, $false, $false, $false, $false, $true, 1, $false, "26+66=", 2) | Out-Null
$d.Content.Find.Execute("48+43=", $true, $false, $false, $false, $false, $true, 1, $false, "45+19=", 2) | Out-Null
$d.Content.Find.Execute("19+13=", $true, $false, $false, $false, $false, $true, 1, $false, "28+25=", 2) | Out-Null
$d.Content.Find.Execute("24+59=", $true, $false, $false, $false, $false, $true, 1, $false, "85-57=", 2) | Out-Null
$d.Content.Find.Execute("28+15=", $true, $false, $false, $false, $false, $true, 1, $false, "66-59=", 2) | Out-Null
$d.Content.Find.Execute("18+56=", $true, $false, $false, $false, $false, $true, 1, $false, "36+39=", 2) | Out-Null
$d.Content.Find.Execute("82-17=", $true, $false, $false, $false, $false, $true, 1, $false, "16+46=", 2) | Out-Null
$d.Content.Find.Execute("6+27=", $true, $false, $false, $false, $false, $true, 1, $false, "4+58=", 2) | Out-Null
$d.Content.Find.Execute("70-64=", $true, $false, $false, $false, $false, $true, 1, $false, "91-74=", 2) | Out-Null
$d.Content.Find.Execute("50-11=", $true, $false, $false, $false, $false, $true, 1, $false, "70-5=", 2) | Out-Null
$d.Content.Find.Execute("44-26=", $true, $false, $false, $false, $false, $true, 1, $false, "38+9=", 2) | Out-Null
$d.Content.Find.Execute("80-76=", $true, $false, $false, $false, $false, $true, 1, $false, "13+58=", 2) | Out-Null
$d.Content.Find.Execute("39+48=", $true, $false, $false, $false, $false, $true, 1, $false, "73-47=", 2) | Out-Null
$d.Content.Find.Execute("9+5=", $true, $false, $false, $false, $false, $true, 1, $false, "22+69=", 2) | Out-Null
$d.Content.Find.Execute("72-65=", $true, $false, $false, $false, $false, $true, 1, $false, "39+47=", 2) | Out-Null
$d.Content.Find.Execute("92-13=", $true, $false, $false, $false, $false, $true, 1, $false, "32-6=", 2) | Out-Null
$d.Content.Find.Execute("90-78=", $true, $false, $false, $false, $false, $true, 1, $false, "61-33=", 2) | Out-Null
$d.Content.Find.Execute("57-49=", $true, $false, $false, $false, $false, $true, 1, $false, "38+5=", 2) | Out-Null
$d.Content.Find.Execute("51-5=", $true, $false, $false, $false, $false, $true, 1, $false, "15+57=", 2) | Out-Null
$d.Content.Find.Execute("15+49=", $true, $false, $false, $false, $false, $true, 1, $false, "19+73=", 2) | Out-Null
$d.Content.Find.Execute("68+19=", $true, $false, $false, $false, $false, $true, 1, $false, "45+8=", 2) | Out-Null
$d.Content.Find.Execute("8+88=", $true, $false, $false, $false, $false, $true, 1, $false, "81-49=", 2) | Out-Null
$d.Content.Find.Execute("69+24=", $true, $false, $false, $false, $false, $true, 1, $false, "55-8=", 2) | Out-Null
$d.Content.Find.Execute("14+59=", $true, $false, $false, $false, $false, $true, 1, $false, "13+39=", 2) | Out-Null
$d.Content.Find.Execute("18+63=", $true, $false, $false, $false, $false, $true, 1, $false, "81-28=", 2) | Out-Null
$d.Content.Find.Execute("8+89=", $true, $false, $false, $false, $false, $true, 1, $false, "77-8=", 2) | Out-Null
$d.Content.Find.Execute("54-29=", $true, $false, $false, $false, $false, $true, 1, $false, "88+9=", 2) | Out-Null
$d.Content.Find.Execute("67+6=", $true, $false, $false, $false, $false, $true, 1, $false, "75-8=", 2) | Out-Null
$d.Content.Find.Execute("49+16=", $true, $false, $false, $false, $false, $true, 1, $false, "36+25=", 2) | Out-Null
$d.Content.Find.Execute("66-17=", $true, $false, $false, $false, $false, $true, 1, $false, "6+68=", 2) | Out-Null
$d.Content.Find.Execute("37+34=", $true, $false, $false, $false, $false, $true, 1, $false, "70-48=", 2) | Out-Null
$d.Content.Find.Execute("91-79=", $true, $false, $false, $false, $false, $true, 1, $false, "63-37=", 2) | Out-Null
$d.Content.Find.Execute("36+57=", $true, $false, $false, $false, $false, $true, 1, $false, "81-52=", 2) | Out-Null
$d.Content.Find.Execute("95-38=", $true, $false, $false, $false, $false, $true, 1, $false, "22+29=", 2) | Out-Null
$d.Content.Find.Execute("58+26=", $true, $false, $false, $false, $false, $true, 1, $false, "89+7=", 2) | Out-Null
$d.Content.Find.Execute("45-18=", $true, $false, $false, $false, $false, $true, 1, $false, "40-9=", 2) | Out-Null
$d.Content.Find.Execute("90-11=", $true, $false, $false, $false, $false, $true, 1, $false, "60-49=", 2) | Out-Null
$d.Content.Find.Execute("9+26=", $true, $false, $false, $false, $false, $true, 1, $false, "28+19=", 2) | Out-Null
$d.Content.Find.Execute("91-49=", $true, $false, $false, $false, $false, $true, 1, $false, "5+48=", 2) | Out-Null
$d.Content.Find.Execute("63+19=", $true, $false, $false, $false, $false, $true, 1, $false, "81-28=", 2) | Out-Null
$d.Content.Find.Execute("54-8=", $true, $false, $false, $false, $false, $true, 1, $false, "9+25=", 2) | Out-Null
$d.Content.Find.Execute("55+37=", $true, $false, $false, $false, $false, $true, 1, $false, "84-15=", 2) | Out-Null
$d.Content.Find.Execute("81-3=", $true, $false, $false, $false, $false, $true, 1, $false, "43-17=", 2) | Out-Null
$d.Content.Find.Execute("41-33=", $true, $false, $false, $false, $false, $true, 1, $false, "54-26=", 2) | Out-Null
$d.Content.Find.Execute("3+18=", $true, $false, $false, $false, $false, $true, 1, $false, "75-36=", 2) | Out-Null
$d.Content.Find.Execute("82-53=", $true, $false, $false, $false, $false, $true, 1, $false, "57+17=", 2) | Out-Null
$d.Content.Find.Execute("54-18=", $true, $false, $false, $false, $false, $true, 1, $false, "51-16=", 2) | Out-Null
$d.Content.Find.Execute("39+36=", $true, $false, $false, $false, $false, $true, 1, $false, "29+43=", 2) | Out-Null
$d.Content.Find.Execute("19+24=", $true, $false, $false, $false, $false, $true, 1, $false, "74-35=", 2) | Out-Null
$d.Content.Find.Execute("90-85=", $true, $false, $false, $false, $false, $true, 1, $false, "95-8=", 2) | Out-Null
$d.Content.Find.Execute("38+18=", $true, $false, $false, $false, $false, $true, 1, $false, "22-17=", 2) | Out-Null
$d.Content.Find.Execute("78-9=", $true, $false, $false, $false, $false, $true, 1, $false, "87-79=", 2) | Out-Null
$d.Content.Find.Execute("27-18=", $true, $false, $false, $false, $false, $true, 1, $false, "66-9=", 2) | Out-Null
$d.Content.Find.Execute("76-37=", $true, $false, $false, $false, $false, $true, 1, $false, "57+24=", 2) | Out-Null
$d.Content.Find.Execute("6+85=", $true, $false, $false, $false, $false, $true, 1, $false, "90-48=", 2) | Out-Null
$d.Content.Find.Execute("18+79=", $true, $false, $false, $false, $false, $true, 1, $false, "65-36=", 2) | Out-Null
$d.Content.Find.Execute("15+47=", $true, $false, $false, $false, $false, $true, 1, $false, "65-37=", 2) | Out-Null
$d.Content.Find.Execute("90-36=", $true, $false, $false, $false, $false, $true, 1, $false, "94-58=", 2) | Out-Null
$d.Content.Find.Execute("91-42=", $true, $false, $false, $false, $false, $true, 1, $false, "25+7=", 2) | Out-Null
$d.Content.Find.Execute("64-48=", $true, $false, $false, $false, $false, $true, 1, $false, "36+38=", 2) | Out-Null
$d.Content.Find.Execute("8+6=", $true, $false, $false, $false, $false, $true, 1, $false, "91-69=", 2) | Out-Null
$d.Content.Find.Execute("48+14=", $true, $false, $false, $false, $false, $true, 1, $false, "37+29=", 2) | Out-Null
$d.Content.Find.Execute("51-12=", $true, $false, $false, $false, $false, $true, 1, $false, "87+8=", 2) | Out-Null
$d.Content.Find.Execute("27+28=", $true, $false, $false, $false, $false, $true, 1, $false, "71-49=", 2) | Out-Null
$d.Content.Find.Execute("84-18=", $true, $false, $false, $false, $false, $true, 1, $false, "44-16=", 2) | Out-Null
$d.Content.Find.Execute("92-55=", $true, $false, $false, $false, $false, $true, 1, $false, "24+48=", 2) | Out-Null
$d.Content.Find.Execute("84-38=", $true, $false, $false, $false, $false, $true, 1, $false, "98-59=", 2) | Out-Null
$d.Content.Find.Execute("23+8=", $true, $false, $false, $false, $false, $true, 1, $false, "57-18=", 2) | Out-Null
$d.Content.Find.Execute("51-44=", $true, $false, $false, $false, $false, $true, 1, $false, "22-16=", 2) | Out-Null
$d.Content.Find.Execute("19+14=", $true, $false, $false, $false, $false, $true, 1, $false, "72-29=", 2) | Out-Null
$d.Content.Find.Execute("92-75=", $true, $false, $false, $false, $false, $true, 1, $false, "22+29=", 2) | Out-Null
$d.Content.Find.Execute("80-28=", $true, $false, $false, $false, $false, $true, 1, $false, "58+6=", 2) | Out-Null

# Step 3: append a brand new row with 5 fresh equations
$newRow = $t.Rows.Add()
$idx = $newRow.Index
$t.Cell($idx,1).Range.Text = "65-17="
$t.Cell($idx,2).Range.Text = "95-8="
$t.Cell($idx,3).Range.Text = "58+29="
$t.Cell($idx,4).Range.Text = "8+19="
$t.Cell($idx,5).Range.Text = "80-1="

Write-Output "done"
